# Refresh the crypto price/volume snapshot (GitHub Actions scheduled update).
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h). Column A holds the static rank index.
#
# A handful of "Price" values are plain decimal numbers (e.g. "17.27", "0.661").
# Assigning those bare through .Value lets Excel auto-detect them as numbers, which
# would change the cell's type away from the source workbook's text cells. We guard
# those with a leading apostrophe -- the normal Excel "store as text" convention --
# so the value round-trips as text exactly like the rest of the column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "36.355.35"
$ws.Range("E2").Value = "  -1.20%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "2.034.95"
$ws.Range("E3").Value = "  -2.66%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  +0.13%  "

# Row 5: BNB
$ws.Range("D5").Value = "'244.73"
$ws.Range("E5").Value = "  -0.65%  "

# Row 6: XRP
$ws.Range("D6").Value = "'0.661"
$ws.Range("E6").Value = "  +1.56%  "

# Row 7: USDC
$ws.Range("E7").Value = "  +0.07%  "

# Row 8: Solana
$ws.Range("D8").Value = "'55.71"
$ws.Range("E8").Value = "  -0.21%  "

# Row 9: OKB
$ws.Range("D9").Value = "'62.40"
$ws.Range("E9").Value = "  +4.89%  "

# Row 10: Cardano
$ws.Range("D10").Value = "'0.364"
$ws.Range("E10").Value = "  -1.11%  "

# Row 11: Dogecoin
$ws.Range("D11").Value = "'0.0742"
$ws.Range("E11").Value = "  -3.35%  "

# Row 12: TRON
$ws.Range("E12").Value = "  -3.02%  "

# Row 13: Polygon
$ws.Range("D13").Value = "'0.897"
$ws.Range("E13").Value = "  +1.85%  "

# Row 14: WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "2.339.95"
$ws.Range("E14").Value = "  -2.31%  "

# Row 15: Chainlink
$ws.Range("E15").Value = "  -5.87%  "

# Row 16: Polkadot
$ws.Range("E16").Value = "  -4.15%  "

# Row 17: WrappedEther
$ws.Range("D17").Value = "2.026.19"
$ws.Range("E17").Value = "  -3.09%  "

# Row 18: WrappedBTC
$ws.Range("D18").Value = "36.255.68"
$ws.Range("E18").Value = "  -1.44%  "

# Row 19: Avalanche
$ws.Range("D19").Value = "'17.27"
$ws.Range("E19").Value = "  -0.98%  "

# Row 20: Litecoin
$ws.Range("E20").Value = "  -2.61%  "

# Row 21: ShibaInu
$ws.Range("D21").Value = "'0.0₃0851"
$ws.Range("E21").Value = "  -3.13%  "

# Row 22: BitcoinCash
$ws.Range("D22").Value = "'236.45"
$ws.Range("E22").Value = "  +0.17%  "

# Row 23: Uniswap
$ws.Range("E23").Value = "  -6.11%  "

# Row 24: Dai
$ws.Range("E24").Value = "  -0.30%  "

# Row 25: Toncoin
$ws.Range("D25").Value = "'2.34"
$ws.Range("E25").Value = "  -2.77%  "

# Row 26: PancakeSwap
$ws.Range("E26").Value = "  +2.21%  "

# Row 27: Cosmos
$ws.Range("D27").Value = "'9.20"
$ws.Range("E27").Value = "  -7.87%  "

# Row 28: Monero
$ws.Range("D28").Value = "'163.24"
$ws.Range("E28").Value = "  -2.99%  "

# Row 29: EthereumClassic
$ws.Range("D29").Value = "'19.88"
$ws.Range("E29").Value = "  -5.64%  "

# Row 30: Stellar
$ws.Range("E30").Value = "  -2.61%  "

# Row 31: ImmutableX
$ws.Range("D31").Value = "'1.20"
$ws.Range("E31").Value = "  -0.62%  "

# Row 32: Filecoin
$ws.Range("D32").Value = "'4.94"
$ws.Range("E32").Value = "  -7.58%  "

# Row 33: Hedera
$ws.Range("D33").Value = "'0.0596"
$ws.Range("E33").Value = "  -2.27%  "

# Row 34: InternetComputer(DFINITY)
$ws.Range("D34").Value = "'4.37"
$ws.Range("E34").Value = "  -7.44%  "

# Row 35: Kaspa
$ws.Range("D35").Value = "'0.0869"
$ws.Range("E35").Value = "  +3.55%  "

# Row 36: BinanceUSD
$ws.Range("E36").Value = "  +0.04%  "

# Row 38: LidoDAOToken
$ws.Range("E38").Value = "  -9.53%  "

# Row 39: THORChain
$ws.Range("E39").Value = "  +1.89%  "

# Row 40: TrustWalletToken
$ws.Range("E40").Value = "  -5.58%  "

# Row 41: HuobiToken
$ws.Range("D41").Value = "'2.86"
$ws.Range("E41").Value = "  -2.14%  "

# Row 42: VeChain
$ws.Range("D42").Value = "'0.0213"
$ws.Range("E42").Value = "  -3.57%  "

# Row 43: ARBITRUM
$ws.Range("E43").Value = "  -5.99%  "

# Row 44: Aave
$ws.Range("D44").Value = "'92.65"
$ws.Range("E44").Value = "  -4.11%  "

# Row 45: Cronos
$ws.Range("D45").Value = "'0.0897"
$ws.Range("E45").Value = "  -6.26%  "

# Row 46: ranking reshuffle swapped Maker and InjectiveProtocol -- row 46 used to be
# Maker and is now InjectiveProtocol.
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").Value = "'15.63"
$ws.Range("E46").Value = "  -4.30%  "

# Row 47: ...and row 47 used to be InjectiveProtocol and is now Maker.
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "1.354.56"
$ws.Range("E47").Value = "  +1.01%  "

# Row 48: FraxShare
$ws.Range("D48").Value = "'7.38"
$ws.Range("E48").Value = "  +4.76%  "

# Row 49: MXToken
$ws.Range("D49").Value = "'2.93"
$ws.Range("E49").Value = "  +1.90%  "

# Row 50: RocketPoolETH
$ws.Range("D50").Value = "2.221.30"
$ws.Range("E50").Value = "  -2.48%  "

# Row 51: MultiversX
$ws.Range("D51").Value = "'45.39"
$ws.Range("E51").Value = "  -0.40%  "
